$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect to make the required edits.
$ws.Unprotect()

# Update the confidential notice text: date changes from 2021-05-19 to 2021-05-20
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2567016281104799
$ws.Range("E2").Value = 0.0005789717461786026

$ws.Range("D3").Value = 0.2551271942098932
$ws.Range("E3").Value = 0.00053777897284224

$ws.Range("D4").Value = 0.243810069027258
$ws.Range("E4").Value = 0.002248069592415192

$ws.Range("D5").Value = 0.2443611086523691
$ws.Range("E5").Value = 0.01714587180163551

$ws.Range("E6").Value = 0.005023711275084386

# Restore sheet protection to its original (protected) state.
$ws.Protect()
